# Update "想去人数" (column F) counts that changed between site generations.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 263
$ws.Range("F3").Value = 612
$ws.Range("F6").Value = 2809
$ws.Range("F9").Value = 29
$ws.Range("F10").Value = 373
$ws.Range("F11").Value = 21
$ws.Range("F12").Value = 311
$ws.Range("F14").Value = 5915
$ws.Range("F16").Value = 1043
$ws.Range("F17").Value = 5
$ws.Range("F18").Value = 100
$ws.Range("F21").Value = 523
$ws.Range("F22").Value = 20
$ws.Range("F23").Value = 19
$ws.Range("F24").Value = 53
$ws.Range("F25").Value = 1294
$ws.Range("F27").Value = 4
$ws.Range("F28").Value = 29
$ws.Range("F29").Value = 2049
$ws.Range("F30").Value = 166
$ws.Range("F31").Value = 345
$ws.Range("F33").Value = 3266

# Sheet 2: 演出 (Performances)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 26
$ws.Range("F7").Value = 343
$ws.Range("F8").Value = 79
$ws.Range("F18").Value = 627
$ws.Range("F20").Value = 61
$ws.Range("F22").Value = 348
$ws.Range("F24").Value = 4046
$ws.Range("F28").Value = 128
$ws.Range("F29").Value = 215
$ws.Range("F30").Value = 68

# Sheet 3: 本地生活 (Local life)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 1802
$ws.Range("F6").Value = 1135
$ws.Range("F8").Value = 1478
$ws.Range("F9").Value = 413
$ws.Range("F12").Value = 629

# Sheet 4: 全部类型 (All types - combined view)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1802
$ws.Range("F6").Value = 1135
$ws.Range("F7").Value = 1478
$ws.Range("F8").Value = 413
$ws.Range("F11").Value = 263
$ws.Range("F12").Value = 612
$ws.Range("F13").Value = 2809
$ws.Range("F15").Value = 29
$ws.Range("F16").Value = 629
$ws.Range("F17").Value = 373
$ws.Range("F18").Value = 79
$ws.Range("F19").Value = 21
$ws.Range("F20").Value = 311
$ws.Range("F22").Value = 5915
$ws.Range("F24").Value = 1043
$ws.Range("F25").Value = 100
$ws.Range("F28").Value = 523
$ws.Range("F32").Value = 61
$ws.Range("F33").Value = 20
$ws.Range("F35").Value = 348
$ws.Range("F36").Value = 1294
$ws.Range("F39").Value = 128
$ws.Range("F40").Value = 215
$ws.Range("F41").Value = 29
$ws.Range("F42").Value = 68
$ws.Range("F44").Value = 2049
$ws.Range("F47").Value = 166
$ws.Range("F48").Value = 345
$ws.Range("F50").Value = 3266
